$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.680.63'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.583.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.06%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("E6").Value = '  -3.19%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.34'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.253'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0590'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0867'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.809.50'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.597.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.531'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.666.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.47%  '
$ws.Range("E20").Value = '  -3.64%  '
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.01%  '
$ws.Range("E29").Value = '  -4.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.40%  '
$ws.Range("E31").Value = '  -3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.377.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.16%  '
$ws.Range("E34").Value = '  -5.17%  '
$ws.Range("E35").Value = '  -5.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.968'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.02%  '
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.537'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.92%  '
$ws.Range("E40").Value = '  -3.80%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  -3.36%  '
$ws.Range("E43").Value = '  -1.84%  '
$ws.Range("E44").Value = '  +1.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.57'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.35%  '
$ws.Range("E46").Value = '  -3.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.719.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0974'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.69%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0497'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.74%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ch1 = [char]0x2087
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0$($ch1)0963"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.34%  '
